# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") values are recalculated/rewritten for rows 2-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 2
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
